$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DetailLogs")
$ws.Activate() | Out-Null

# Round-off time entry values corrected for the Detail Logs test data
$ws.Range("B2").Value = "4.777"
$ws.Range("D2").Value = "4.8"

# Leave the cursor where the author left it before saving
$ws.Range("D10").Select() | Out-Null
